# Update to "Denmark Division 3" sheet: several match rows had their data
# shuffled between row positions (rows B..AB, i.e. everything except the
# row-index column A). Column A keeps the sequential row id, so only the
# match data (id, HomeTeam/AwayTeam, odds, PL columns, ...) moves.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Swap-Rows($rowA, $rowB) {
    $rangeA = $ws.Range("B$rowA`:AB$rowA")
    $rangeB = $ws.Range("B$rowB`:AB$rowB")
    $valA = $rangeA.Value()
    $valB = $rangeB.Value()
    $rangeA.Value = $valB
    $rangeB.Value = $valA
}

# Simple pairwise swaps of row data
Swap-Rows 62 64
Swap-Rows 66 67
Swap-Rows 70 71
Swap-Rows 162 163
Swap-Rows 164 166

# 3-way rotation for rows 146, 147, 148:
# new146 = old148, new147 = old146, new148 = old147
$r146 = $ws.Range("B146:AB146")
$r147 = $ws.Range("B147:AB147")
$r148 = $ws.Range("B148:AB148")
$v146 = $r146.Value()
$v147 = $r147.Value()
$v148 = $r148.Value()
$r146.Value = $v148
$r147.Value = $v146
$r148.Value = $v147
